$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values per repulled data
$ws.Range("F4").Value = 4
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = 0
$ws.Range("F7").Value = -4
$ws.Range("F8").Value = 1
$ws.Range("F11").Value = -3
$ws.Range("F13").Value = 0
$ws.Range("F15").Value = 2
$ws.Range("F16").Value = 2
$ws.Range("F17").Value = 3
$ws.Range("F19").Value = -6
